$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix email values that had the number appended instead of inserted before @mail.com
$ws.Range("G4").Value = "nemreg1es1@mail.com"
$ws.Range("G5").Value = "nemreg1es2@mail.com"
$ws.Range("G6").Value = "nemreg1es3@mail.com"

# Remove the now-duplicated / unwanted email entries entirely
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Update the active selection to reflect where the user ended up (G4)
$ws.Range("G4").Select()
